$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Validation Phase results (columns B, C, D) for rows 3-12
$ws.Range("B3").Value = 0.64444195089672496
$ws.Range("C3").Value = 0.57679420009771498
$ws.Range("D3").Value = 0.0012537402340823799

$ws.Range("B4").Value = 0.639741755592796
$ws.Range("C4").Value = 0.66737576791656394
$ws.Range("D4").Value = 0.00126594542567377

$ws.Range("B5").Value = 0.64400234772562404
$ws.Range("C5").Value = 0.58001209282675303
$ws.Range("D5").Value = 0.00125437175499044

$ws.Range("B6").Value = 0.75932645797729403
$ws.Range("C6").Value = 6.8614125000000001
$ws.Range("D6").Value = 0.0048430510000000001

$ws.Range("B7").Value = 0.85400037812962104
$ws.Range("C7").Value = 0.428416795867057
$ws.Range("D7").Value = 0.00072027827221347201

$ws.Range("B8").Value = 0.79095507517726504
$ws.Range("C8").Value = 0.53355472755182398
$ws.Range("D8").Value = 0.00093256825805548003

$ws.Range("B9").Value = 0.90717716496461898
$ws.Range("C9").Value = 0.33961996526183902
$ws.Range("D9").Value = 0.00047667283584392801

$ws.Range("B10").Value = 0.96639268883587004
$ws.Range("C10").Value = 0.267330764018261
$ws.Range("D10").Value = 0.000230317100128948

$ws.Range("B11").Value = 0.97038525885994098
$ws.Range("C11").Value = 0.28070083835049903
$ws.Range("D11").Value = 0.00020328828731852101

$ws.Range("B12").Value = 0.80308792299610798
$ws.Range("C12").Value = 0.99511655463733995
$ws.Range("D12").Value = 0.00083865796941974698

# Update the active cell selection to reflect the cursor position after data entry
[void]$ws.Range("B14").Select()
